$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: new work entry (date, start time, end time, description)
$ws.Range("A31").Value = 44096
$ws.Range("A31").NumberFormat = $ws.Range("A30").NumberFormat

$ws.Range("B31").Value = 0.66666666666666663
$ws.Range("B31").NumberFormat = $ws.Range("B30").NumberFormat

$ws.Range("C31").Value = 0.91666666666666663
$ws.Range("C31").NumberFormat = $ws.Range("C30").NumberFormat

$ws.Range("E31").Value = "Bugfixes + Frontpage + Randomgen"

# Update selection to reflect where the author left off editing
$ws.Range("E31").Select() | Out-Null
